# Update "想去人数" (want-to-go count) figures that changed between scrapes.
$wb = $excel.ActiveWorkbook

# Sheet "展览": row 4 -> F4 1494 -> 1499, row 5 -> F5 697 -> 698
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 1499
$wsExhibit.Range("F5").Value = 698

# Sheet "全部类型": row 4 -> F4 1494 -> 1499, row 6 -> F6 697 -> 698
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1499
$wsAll.Range("F6").Value = 698
